# Update Civil Works Expenditure figures across the four sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "completed" ---
$ws = $wb.Worksheets.Item("completed")
$ws.Range("D2").Value = 599
$ws.Range("L2").Value = 599

# --- Sheet "remaining" ---
$ws = $wb.Worksheets.Item("remaining")
$ws.Range("D2").Value = 599
$ws.Range("E2").Value = 599
$ws.Range("F2").Value = 599
$ws.Range("G2").Value = 589.88
$ws.Range("H2").Value = 524.25
$ws.Range("I2").Value = 427.21
$ws.Range("J2").Value = 329.7
$ws.Range("K2").Value = 169.7

# --- Sheet "Current_month" ---
$ws = $wb.Worksheets.Item("Current_month")
$ws.Range("D2").Value = 599
$ws.Range("L2").Value = 169.7

# --- Sheet "Prev_cum" ---
$ws = $wb.Worksheets.Item("Prev_cum")
$ws.Range("D2").Value = 599
